$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    # Force a literal (non-numeric-coerced) text value into a cell: build it as a
    # formula returning the literal string, then paste-special as values only so
    # the stored cell is plain text with no residual formula and no style change.
    $rng.Formula = $val
    $rng.Copy() | Out-Null
    $rng.PasteSpecial(-4163) | Out-Null
}

Set-TextValue $ws.Range("D2") '="26.268.83"'
$ws.Range("E2").Value = '  -0.16%  '

Set-TextValue $ws.Range("D3") '="1.594.47"'
$ws.Range("E3").Value = '  +0.21%  '

$ws.Range("E4").Value = '  -0.05%  '

Set-TextValue $ws.Range("D5") '="213.12"'
$ws.Range("E5").Value = '  +0.55%  '

Set-TextValue $ws.Range("D6") '="0.498"'
$ws.Range("E6").Value = '  -0.52%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -0.28%  '

$ws.Range("E9").Value = '  -0.43%  '

Set-TextValue $ws.Range("D10") '="18.99"'
$ws.Range("E10").Value = '  -1.88%  '

Set-TextValue $ws.Range("D11") '="0.0849"'
$ws.Range("E11").Value = '  +0.14%  '

Set-TextValue $ws.Range("D12") '="1.818.67"'
$ws.Range("E12").Value = '  +0.21%  '

Set-TextValue $ws.Range("D13") '="1.582.06"'
$ws.Range("E13").Value = '  -1.36%  '

$ws.Range("E14").Value = '  -1.09%  '

$ws.Range("E15").Value = '  -2.11%  '

Set-TextValue $ws.Range("D16") '="63.84"'
$ws.Range("E16").Value = '  -0.97%  '

Set-TextValue $ws.Range("D17") '="26.256.27"'
$ws.Range("E17").Value = '  -0.24%  '

Set-TextValue $ws.Range("D18") '="0.0₃0722"'
$ws.Range("E18").Value = '  -1.32%  '

Set-TextValue $ws.Range("D19") '="215.16"'
$ws.Range("E19").Value = '  +1.46%  '

Set-TextValue $ws.Range("D20") '="7.35"'
$ws.Range("E20").Value = '  -1.67%  '

$ws.Range("E21").Value = '  -0.05%  '

$ws.Range("E22").Value = '  -0.05%  '

Set-TextValue $ws.Range("D23") '="9.03"'
$ws.Range("E23").Value = '  +0.29%  '

$ws.Range("E24").Value = '  -2.68%  '

Set-TextValue $ws.Range("D25") '="144.80"'
$ws.Range("E25").Value = '  -0.31%  '

$ws.Range("E26").Value = '  -0.04%  '

Set-TextValue $ws.Range("D27") '="6.97"'
$ws.Range("E27").Value = '  -1.18%  '

$ws.Range("E28").Value = '  +0.86%  '

$ws.Range("E29").Value = '  -0.58%  '

$ws.Range("E30").Value = '  -2.12%  '

$ws.Range("E32").Value = '  -0.64%  '

Set-TextValue $ws.Range("D33") '="1.417.57"'
$ws.Range("E33").Value = '  +5.90%  '

$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("E35").Value = '  -0.77%  '

$ws.Range("E36").Value = '  -1.62%  '

Set-TextValue $ws.Range("D37") '="0.573"'
$ws.Range("E37").Value = '  -4.93%  '

$ws.Range("E38").Value = '  -0.97%  '

Set-TextValue $ws.Range("D39") '="0.823"'
$ws.Range("E39").Value = '  +0.60%  '

$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("E41").Value = '  -0.04%  '

$ws.Range("E42").Value = '  +0.83%  '

Set-TextValue $ws.Range("D43") '="0.925"'
$ws.Range("E43").Value = '  -12.42%  '

$ws.Range("E44").Value = '  -0.29%  '

Set-TextValue $ws.Range("D45") '="1.730.24"'
$ws.Range("E45").Value = '  +0.18%  '

Set-TextValue $ws.Range("D46") '="60.93"'

Set-TextValue $ws.Range("D47") '="87.54"'
$ws.Range("E47").Value = '  -0.40%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D48") '="0.0₆0104"'
$ws.Range("E48").Value = '  -0.09%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D49") '="1.48"'
$ws.Range("E49").Value = '  -1.61%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D50") '="0.0502"'
$ws.Range("E50").Value = '  -0.38%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D51") '="0.0954"'
$ws.Range("E51").Value = '  -2.74%  '
